# Performs tare after calibration.
# Adds a "Célula de carga" (load cell) line item to the budget sheet and
# clarifies several component labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new row above the old "Placa ilhada" row for the new
# "Célula de carga" line item; this shifts rows 5:11 down to 6:12 and the
# SUM() formulas in the (now) Total row auto-extend their ranges.
$ws.Rows.Item(5).Insert()

# The freshly inserted row inherits default formatting; copy the
# (identical) formats from the row below it so row 5 matches the rest of
# the table (currency format + borders on B/C, bordered text on A).
$ws.Range("A6:C6").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row and rename a few labels (order chosen so that new
# shared-string entries are appended in the same sequence as the source
# workbook).
$ws.Range("A5").Value = "Célula de carga"
$ws.Range("A9").Value = "Cabo USB-USB-tipoB (por célula de carga)"
$ws.Range("A6").Value = "Placa ilhada (7cm x 9cm)"
$ws.Range("A7").Value = "Conector USB tipo B"

# Column A needs to be a bit wider to fit the longer labels.
$ws.Columns.Item(1).ColumnWidth = 34.92

# Match the author's final selection in the saved file.
$ws.Activate() | Out-Null
$ws.Range("A3:A9").Select() | Out-Null
